$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data for these cells was regenerated with extra decimal
# precision. The cells hold their numbers as literal text (shared strings,
# not numeric cells), so each target cell is first forced to the "Text"
# number format before the new value is assigned -- this keeps Excel's
# automatic "looks like a number" conversion from turning the literal back
# into a numeric cell type.

# Enterprises (absolute #): Micro 385826.5 -> 385826.54, SMEs 14408.5 -> 14408.46
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "385826.54"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "14408.46"

# Enterprises density (per 1000 people): Micro 19.7 -> 19.67, SMEs 0.7 -> 0.73, MSMEs 20.4 -> 20.41
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "19.67"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "0.73"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.41"

# Employment (% of total): MSMEs 26.4 -> 26.39
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.39"
